$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Staatssteuer")
$ws2 = $wb.Worksheets.Item("Bundessteuer")

# --- Add new tariff type "single" (Tariftyp = 2) to Staatssteuer ---
# The existing 12 data rows (2-13, Tariftyp=1) are copied down to rows 14-25
# to seed the new tariff type, then the Tariftyp is switched to 2 there.
$ws1.Range("A2:H13").Copy()
$ws1.Range("A14").PasteSpecial()
$excel.CutCopyMode = $false
$ws1.Range("D14:D25").Value = 2

# The original tariff (Tariftyp=1, rows 2-13) gets updated Einkommen/Steuer values.
$updatedValues = @(
    @(6700, 0),
    @(11400, 94),
    @(16100, 235),
    @(23700, 539),
    @(33000, 1004),
    @(43700, 1646),
    @(56100, 2514),
    @(73000, 3866),
    @(105500, 6791),
    @(137700, 10011),
    @(188700, 15621),
    @(254900, 23565)
)
for ($i = 0; $i -lt $updatedValues.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 5).Value = $updatedValues[$i][0]
    $ws1.Cells.Item($row, 6).Value = $updatedValues[$i][1]
}

# --- Update view/selection state ---
# Bundessteuer is no longer the active tab / no longer scrolled; its selection moves.
$ws2.Activate()
$ws2.Range("I20").Select()

# Staatssteuer becomes the active tab, with its own new selection.
$ws1.Activate()
$ws1.Range("F12").Select()
